$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.086.77"
$ws.Range("E2").Value = "  -0.40%  "

$ws.Range("D3").Value = "1.875.77"
$ws.Range("E3").Value = "  -1.85%  "

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.003"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.11%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "313.23"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.5039"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.47%  "

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3839"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  -2.19%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.08557"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -7.78%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "1.116"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  -2.22%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "41.28"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -1.94%  "

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "6.289"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.89%  "

$ws.Range("E13").Value = "  -1.62%  "

$ws.Range("D14").Value = "1.859.24"
$ws.Range("E14").Value = "  -2.38%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "7.211"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -1.62%  "

$ws.Range("E16").Value = "  +0.16%  "

$ws.Range("E17").Value = "  -2.16%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "91.08"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -1.72%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "0.06629"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  +0.05%  "

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "18.04"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "  +0.17%  "

$ws.Range("E21").Value = "  +0.17%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.090"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -2.29%  "

$ws.Range("D23").Value = "28.126.53"

$ws.Range("E24").Value = "  -1.67%  "

$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "2.269"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.52%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.595"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.08%  "

$ws.Range("D27").Value = "2.086.14"
$ws.Range("E27").Value = "  -2.01%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "20.73"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.93%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "156.76"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  -0.94%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "126.29"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.87%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.1059"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.36%  "

$ws.Range("E32").Value = "  -4.04%  "

$ws.Range("E33").Value = "  -1.17%  "

$ws.Range("E34").Value = "  -0.54%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.630"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.70%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.02460"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  +0.87%  "

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.06564"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -2.05%  "

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "0.2184"
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = "  -2.10%  "

$ws.Range("E39").Value = "  -2.35%  "

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "1.240"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  -3.53%  "

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.6397"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -2.77%  "

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "11.49"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.64%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "4.895"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  -2.47%  "

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "13.25"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  -0.69%  "

$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.6016"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  -2.42%  "

$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.283"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.14%  "

$ws.Range("E47").Value = "  -1.51%  "

$ws.Range("E48").Value = "  -1.84%  "

$ws.Range("E49").Value = "  +2.61%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "121.32"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -0.68%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "80.51"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +2.82%  "
